$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tests")

# Two new experiment rows are appended to the "Tests" sheet (row 35 is left
# blank as a separator, matching the existing layout of the table).
#
# The "Same" string must be registered in the shared-string table before
# "Basic set-up, ..." so that new shared-string indices line up with the
# target workbook, so row 37's text is written first.
$ws.Cells.Item(37, 1).Value = "Wessel1"
$ws.Cells.Item(37, 9).Value = "Same"

# Row 36: Wessel1 / Standardtest2 / w approach
$ws.Cells.Item(36, 1).Value = "Wessel1"
$ws.Cells.Item(36, 2).Value = "Standardtest2"
$ws.Cells.Item(36, 3).Value = 3
$ws.Cells.Item(36, 4).Value = 3
$ws.Cells.Item(36, 5).Value = "w"
$ws.Cells.Item(36, 6).Value = 15
$ws.Cells.Item(36, 7).Value = 20
$ws.Cells.Item(36, 8).Value = 13
$ws.Cells.Item(36, 9).Value = "Basic set-up, but timblserver now and new calculation of cks and skks… but possibly too slow nl corpus"

# Finish row 37: Wessel1 / 10% attenuation / w approach
$ws.Cells.Item(37, 2).Value = 0.1
# Reuse the existing percentage number format (as already used by B23/B26/.../B34)
# instead of applying a named style, so no new style gets introduced.
$ws.Range("B37").NumberFormat = $ws.Range("B34").NumberFormat
$ws.Cells.Item(37, 3).Value = 3
$ws.Cells.Item(37, 4).Value = 3
$ws.Cells.Item(37, 5).Value = "w"
$ws.Cells.Item(37, 6).Value = 33
$ws.Cells.Item(37, 7).Value = 41
$ws.Cells.Item(37, 8).Value = 2661

$ws.Range("I37").Select()
